$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81, shifting existing rows 81..108 down to 82..109
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new weekly record
$ws.Cells.Item(81, 1).Value = 3
$ws.Cells.Item(81, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(81, 3).Value = "Coquimbo"
$ws.Cells.Item(81, 4).Value = 44559
$ws.Cells.Item(81, 5).Value = 5
$ws.Cells.Item(81, 6).Value = 100112052
$ws.Cells.Item(81, 7).Value = "Albahaca"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 100
$ws.Cells.Item(81, 11).Value = 4000
$ws.Cells.Item(81, 12).Value = 4500
$ws.Cells.Item(81, 13).Value = 4250
$ws.Cells.Item(81, 14).Value = "$/docena de matas"
$ws.Cells.Item(81, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(81, 16).Value = 708
$ws.Cells.Item(81, 17).Value = 6
$ws.Cells.Item(81, 18).Value = "Hortaliza"
